$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-28 Saturday" "2024-12-29 Sunday"

Replace-Text "459×2=" "692×6="
Replace-Text "512×5=" "413×8="
Replace-Text "151×3=" "369×5="
Replace-Text "942×3=" "765×6="
Replace-Text "921×5=" "427×6="
Replace-Text "218×6=" "153×2="
Replace-Text "968×4=" "297×9="
Replace-Text "385×7=" "728×4="
Replace-Text "598×6=" "712×7="
Replace-Text "307×6=" "885×9="
Replace-Text "293×3=" "446×3="
Replace-Text "245×7=" "652×6="
Replace-Text "162×5=" "803×7="
Replace-Text "866×6=" "151×9="
Replace-Text "255×8=" "731×8="
Replace-Text "390×2=" "209×9="
Replace-Text "294×9=" "438×8="
Replace-Text "424×6=" "571×9="
Replace-Text "761×6=" "320×2="
Replace-Text "946×3=" "898×5="
Replace-Text "436×4=" "578×2="
Replace-Text "179×4=" "491×7="
Replace-Text "699×6=" "234×2="
Replace-Text "742×2=" "646×9="
Replace-Text "890×7=" "818×5="
